# Update the Yearly sheet (401K and Suzie's Roth IRA columns for the
# September 2017 row), then update the selection/active cell on both
# sheets to match the author's final cursor position.

$wb = $excel.ActiveWorkbook

$yearly = $wb.Worksheets.Item("Yearly")
$allTime = $wb.Worksheets.Item("All Time")

# Row 11 = September (month 9) under the 2017 block (columns J:O).
# M11 = 401K dividends, N11 = Suzie's Roth IRA dividends.
# O11 = SUM(L11:N11) recalculates automatically.
$yearly.Range("M11").Value = 30.05
$yearly.Range("N11").Value = 12.11

# Recalculate so dependent formulas (O11, M15:O15, and the cross-sheet
# formulas on "All Time") pick up the new totals.
$excel.Calculate()

# Update the saved cursor/selection position on each sheet.
$yearly.Activate()
$yearly.Range("G26").Select()

$allTime.Activate()
$allTime.Range("O32").Select()

$excel.Calculate()
